$wb = $excel.ActiveWorkbook
$ws2 = $wb.Worksheets.Item("Debits")
$ws2.Columns.Item(1).ColumnWidth = 15.9
Write-Host "done"
